# Fix 25: address issues handled
# Fills in the missing "address" values (column D) for rows where the
# address text had previously only been duplicated into the description
# column (H), leaving the dedicated address column blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$addresses = @{
    4  = "Redwood, New York"
    5  = "Durham, North Carolina"
    6  = "New Milford, Connecticut"
    9  = "Mackinaw City, Michigan"
    10 = "Gilboa, New York`nStamford, New York"
    11 = "Oklahoma City, Oklahoma"
    12 = "La Riviere, Manitoba, Canada"
    13 = "Vestal, New York"
    14 = "Broken Bow, Oklahoma"
    15 = "Al Marina, Abu Dhabi, United Arab Emirates"
    18 = "San Diego, California"
    19 = "St. Petersberg, Russia"
}

foreach ($row in $addresses.Keys) {
    $ws.Range("D$row").Value = $addresses[$row]
}

# Row 10's address contains an embedded line break (same text already shown
# in H10); writing it into D10 must not change the row's rendered height.
$ws.Rows.Item(10).RowHeight = 12.75

# Widen column D to fit the newly-populated addresses (target stored width
# is 37.7142857142857 character-units; 36.83 is the closest input that the
# runtime's column-width rounding maps onto that value).
$ws.Columns.Item(4).ColumnWidth = 36.83
